$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case field names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the "de"/"del"/"el"/"la" connectors in place names ---
$ws.Range("B13").Value = "Hidalgo Del Parral"
$ws.Range("B15").Value = "Valle De Zaragoza"
$ws.Range("A17").Value = "Ciudad De México"
$ws.Range("A27").Value = "Estado De México"
$ws.Range("B27").Value = "Almoloya De Alquisiras"
$ws.Range("B28").Value = "Ecatepec De Morelos"
$ws.Range("B37").Value = "San Miguel De Allende"
$ws.Range("B40").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B46").Value = "Silao De La Victoria"
$ws.Range("B48").Value = "Coyuca De Catalán"
$ws.Range("B53").Value = "Tlapa De Comonfort"
$ws.Range("B56").Value = "Mineral Del Chico"
$ws.Range("B58").Value = "Tenango De Doria"
$ws.Range("B59").Value = "Tulancingo De Bravo"
$ws.Range("B77").Value = "Putla Villa De Guerrero"
$ws.Range("B78").Value = "San José Del Progreso"
$ws.Range("B80").Value = "San Miguel El Grande"
$ws.Range("B90").Value = "Tlacolula De Matamoros"
$ws.Range("B91").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B94").Value = "Palmar De Bravo"
$ws.Range("B95").Value = "San Salvador El Verde"
$ws.Range("B99").Value = "Pinal De Amoles"
$ws.Range("B116").Value = "Lerdo De Tejada"

# --- Remove trailing footer/metadata rows (127-131), shrinking the used range to A1:D125 ---
$ws.Rows("127:131").Delete()
